$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date from 2023-09-01 to 2023-09-05
# for every data row (rows 2-12).
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = "2023-09-05"
}
